$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (Product Backlog column titles) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Story"
$ws.Range("C1").Value = "Sprint"
$ws.Range("D1").Value = "Story Priority"
$ws.Range("E1").Value = "Story Status"
$ws.Range("F1").Value = "Story Points"

# Apply the built-in "20% - Accent3" cell style to the header row
$ws.Range("A1:F1").Style = "20% - Accent3"

# --- ID column sample data (rows 2-12 => 0 through 10) ---
for ($i = 0; $i -le 10; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# --- Column widths for the backlog fields ---
$ws.Columns.Item(2).ColumnWidth = 24.9167
$ws.Columns.Item(3).ColumnWidth = 21.9167
$ws.Columns.Item(4).ColumnWidth = 25.4167
$ws.Columns.Item(5).ColumnWidth = 26.75
$ws.Columns.Item(6).ColumnWidth = 22.9167

# --- Selection / active cell ---
[void]$ws.Range("B5").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
